$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 2 (Admin/admin123), shifting it down to row 4
$ws.Rows.Item(2).Resize(2).Insert()

# New row 2: admin / ad123
$ws.Range("A2").Value = "admin"
$ws.Range("B2").Value = "ad123"

# New row 3: ad123 / admin
$ws.Range("A3").Value = "ad123"
$ws.Range("B3").Value = "admin"

# Set the active cell / selection to match the target state
$ws.Range("B8").Select()
